# Add "accessType" column to the "queue" sheet, update the related named
# range, and move the active-sheet/selection state around (as captured by
# the workbook's cached view state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "queue" worksheet: insert a new "accessType" column after "enabled"
# ---------------------------------------------------------------------
$queueSheet = $wb.Worksheets.Item("queue")

# Insert a new column C (pushes former C..I to D..J)
$queueSheet.Columns("C").Insert()
$queueSheet.Columns("C").ColumnWidth = 9.5

# Header
$queueSheet.Cells.Item(1, 3).Value = "accessType"

# Data rows
$queueSheet.Cells.Item(2, 3).Value = "exclusive"
$queueSheet.Cells.Item(3, 3).Value = "non-exclusive"
$queueSheet.Cells.Item(4, 3).Value = "non-exclusive"

# Update the "queue" named range to account for the extra column
foreach ($n in $wb.Names) {
    if ($n.Name -eq "queue!queue") {
        $n.RefersTo = "=queue!`$A`$1:`$J`$4"
    }
}

# Restore the (now stale) selection on the queue sheet
$queueSheet.Range("C9").Select()

# ---------------------------------------------------------------------
# 2. "bridge" worksheet: no longer the active tab, selection unchanged
# ---------------------------------------------------------------------
$bridgeSheet = $wb.Worksheets.Item("bridge")
$bridgeSheet.Range("D3").Select()

# ---------------------------------------------------------------------
# 3. "msg-vpn" worksheet: becomes the active tab/sheet
# ---------------------------------------------------------------------
$msgVpnSheet = $wb.Worksheets.Item("msg-vpn")
$msgVpnSheet.Activate()
$msgVpnSheet.Range("A2").Select()
